$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos (PT) text was missing; fill it in ---
$ws.Range("B10:C10").Value = 'Propiciar ao discente conhecimento dos fundamentos da Educação Ambiental utilizando como base os problemas ambientais da atualidade. Desenvolver atividades práticas integradas à região. Orientar o desenvolvimento de projetos relacionados à Gestão e Educação Ambiental.'

# --- Restructure rows 13-21 into the new row 13-23 layout ---
# Remove the old (misaligned) rows entirely so no stray formatting remains
$ws.Range("A13:A21").EntireRow.Delete()
# Insert 11 fresh blank rows in their place (rows 13-23)
$ws.Range("A13:A23").EntireRow.Insert()

# Re-apply the B/C column formatting (wrap text, font) used throughout the sheet
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C23").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
# The row insert + paste leaves stray empty styled cells in column A; remove them
$ws.Range("A13:A14").Clear()

# Row 13: Docentes responsáveis (1st professor) - no label in col A
$ws.Range("B13:C13").Value = '9146830 - Danúbia Caporusso Bargos'

# Row 14: Docentes responsáveis (2nd professor) - no label in col A
$ws.Range("B14:C14").Value = '5817650 - Érica Leonor Romão'

# Row 15: Programa resumido:
$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15:C15").Value = 'Considerações gerais sobre a problemática ambiental. Evolução das questões ambientais no Brasil e no mundo. Educação e Gestão Ambiental. Elaboração e acompanhamento de projetos de educação ambiental.'
$ws.Rows.Item(15).RowHeight = 60

# Row 16: Short syllabus:
$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16:C16").Value = 'General considerations on environmental problem. Evolution of environmental questions in Brazil and in the world. Education and Environmental Management. Development and monitoring of environmental education projects.'
$ws.Rows.Item(16).RowHeight = 60

# Row 17: Programa:
$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17:C17").Value = 'Sociedade, natureza e desenvolvimento. A relação degradação ambiental-qualidade de vida. Meio ambiente e cidadania. Percepção e Interpretação ambiental. Meio ambiente e representação social. Histórico da educação ambiental e conceitos de meio ambiente; Conceitos, princípios e pensamentos norteadores da Educação Ambiental. A questão ambiental e as conferências mundiais de meio ambiente. O movimento ambientalista e o histórico da EA no Brasil e no mundo; A Agenda 21 e educação ambiental. A política nacional de educação ambiental (pnea) e legislação correlata: A abordagem interdisciplinar da educação ambiental; Educação como instrumento de Gestão Ambiental. Educação ambiental nas empresas e o Sistema de Gestão Ambiental. Projetos, reflexões e práticas da Educação Ambiental. Análise e vivências de experiências práticas de educação ambiental em diferentes contextos. Metodologia de projetos, oficinas e capacitação em educação ambiental.'
$ws.Rows.Item(17).RowHeight = 120

# Row 18: Syllabus:
$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18:C18").Value = 'Society, nature and development. The relation environmental degradation-quality of life. Environment and citizenship. Perception and Environmental Interpretation. Environment and social representation. History of environmental education and environmental concepts; Concepts, principles and thoughts guiding Environmental Education. The environmental issue and the world environmental conferences. The environmental movement and EA''s history in Brazil and in the world; Agenda 21 and environmental education. The national policy of environmental education (pnea) and related legislation: The interdisciplinary approach of environmental education; Education as an instrument of Environmental Management. Environmental education in companies and the Environmental Management System. Projects, reflections and practices of Environmental Education. Analysis and experiences of practical experiences of environmental education in different contexts. Methodology of projects, workshops and training in environmental education'
$ws.Rows.Item(18).RowHeight = 120

# Row 19: Avaliação: (label only, no B/C content)
$ws.Range("A19").Value = 'Avaliação:'
$ws.Range("B19:C19").Clear()

# Row 20: Método:
$ws.Range("A20").Value = 'Método:'
$ws.Range("B20:C20").Value = 'Avaliação baseada em provas, exercícios, projetos, seminários e outras formas de avaliação, sendo a nota final correspondente a média ponderada das notas atribuídas às avaliações aplicadas'
$ws.Rows.Item(20).RowHeight = 60

# Row 21: Critério:
$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21:C21").Value = 'Nota Final: NF ≥ 5,0'
$ws.Rows.Item(21).RowHeight = 60

# Row 22: Norma de recuperação:
$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22:C22").Value = 'Provas e/ou exercícios dirigidos'
$ws.Rows.Item(22).RowHeight = 60

# Row 23: Bibliografia:
$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B23:C23").Value = 'CARVALHO, I. C. M.; Educação Ambiental e formação do sujeito ecológico. São Paulo: Cortez, 2006.CINQUETTI, H. C. S.; LOGAREZZI, A. (Org.). Consumo e Resíduo - Fundamentos para o trabalho educativo. 1 ed. São Carlos: EdUFSCar, 2006, v. 1.DIAS, G. F. Dinâmica e instrumentação para educação ambiental. 1. ed. São Paulo: Gaia, 2010. v. 1. 216p.DIAS, G. F. Educação e Gestão Ambiental. 1. ed. São Paulo: Editora Gaia Ltda, 2006. v. 1. 118p.DIAS, G. F. Educação Ambiental: princípios e práticas. 6a ed. São Paulo: Gaia, 2000.GUIMARÃES, M. (org.) Caminhos da educação ambiental: da forma à ação. Campinas, SP: Papirus, 2006.JACOBI, Pedro Roberto, MONTEIRO,F. M ; FERNANDES, M. L. B. . Educação e Sustentabilidade- caminhos e práticas para uma educação transformadora. São Paulo: Evoluir Cultural, 2009. v. 01. 108p.JACOBI, Pedro Roberto OLIVEIRA, F. C. J. F. (Org.). Educação, Meio Ambiente e Cidadania - reflexões e experiências. São Paulo: SMA/CEAM, 1998. 121p LOUREIRO, C. F. B. Trajetória e fundamentos da educação ambiental. 4. ed. São Paulo: Cortez editora, 2012. 165pPHILIPPI JR., A & PELICIONI, M. C. F. (Eds). 2005. Educação ambiental e sustentabilidade. Barueri SP: Manole. 878p. (Coleção Ambiental, 3).'
$ws.Rows.Item(23).RowHeight = 120

# --- Column layout cleanup: column A should only be 30.71 wide on its own ---
# (column B keeps its existing 60.71 width/style)
$ws.Columns.Item(2).ColumnWidth = 59.877604166666668

